$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Forum Orang Tua Siswa Minta MBG di Jabar Disetop Sementara', '2025-09-29T13:32:18+07:00', 'Bima Bagaskara', 'https://news.detik.com/berita/d-8135596/forum-orang-tua-siswa-minta-mbg-di-jabar-disetop-sementara')
    ,@('Belasan Siswa MTS di KBB Kembali Dirawat Usai Gejala Keracunan MBG Kambuh Lagi', '2025-09-29T13:16:52+07:00', 'Whisnu Pradana', 'https://news.detik.com/berita/d-8135570/belasan-siswa-mts-di-kbb-kembali-dirawat-usai-gejala-keracunan-mbg-kambuh-lagi')
    ,@('Harga Rp 10 Ribu, Ini Penampakan Menu Dapur Sehat SD di Solo yang Tolak MBG', '2025-09-29T13:12:02+07:00', 'Tara Wahyu NV', 'https://www.detik.com/jateng/berita/d-8135560/harga-rp-10-ribu-ini-penampakan-menu-dapur-sehat-sd-di-solo-yang-tolak-mbg')
    ,@('Sehari 12 Jam di Dapur, Chef Program MBG Ini Ungkap Kesibukannya', '2025-09-29T13:00:04+07:00', 'Diah Afrilian', 'https://food.detik.com/info-kuliner/d-8135505/sehari-12-jam-di-dapur-chef-program-mbg-ini-ungkap-kesibukannya')
    ,@('Belasan Siswa MTs di Cipongkor Alami Gejala Berulang Keracunan MBG', '2025-09-29T12:46:07+07:00', 'Whisnu Pradana', 'https://www.detik.com/jabar/berita/d-8135490/belasan-siswa-mts-di-cipongkor-alami-gejala-berulang-keracunan-mbg')
    ,@('Bingkai Sepekan: MBG Jadi Polemik, Menu Minim Gizi-Ribuan Keracunan', '2025-09-27T07:00:56+07:00', '-', 'https://health.detik.com/fotohealth/d-8132256/bingkai-sepekan-mbg-jadi-polemik-menu-minim-gizi-ribuan-keracunan')
    ,@('''''Tur'' ke Dapur Makan Bergizi Gratis di Bandung', '2025-09-26T19:00:53+07:00', '-', 'https://www.detik.com/jabar/foto/d-8131846/tur-ke-dapur-makan-bergizi-gratis-di-bandung')
    ,@('MBG di Sekolah Jakut Isinya Cuma Snack, Gizi Dipertanyakan', '2025-09-26T17:00:17+07:00', '-', 'https://health.detik.com/fotohealth/d-8131726/mbg-di-sekolah-jakut-isinya-cuma-snack-gizi-dipertanyakan')
    ,@('Penampakan Dapur MBG di Bandung yang Bikin Ratusan Siswa Keracunan', '2025-09-25T18:00:53+07:00', '-', 'https://news.detik.com/foto-news/d-8130038/penampakan-dapur-mbg-di-bandung-yang-bikin-ratusan-siswa-keracunan')
    ,@('Potret Pilu Ribuan Anak Sekolah Jadi Korban Makan Bergizi Gratis', '2025-09-25T13:10:51+07:00', '-', 'https://health.detik.com/fotohealth/d-8129561/potret-pilu-ribuan-anak-sekolah-jadi-korban-makan-bergizi-gratis')
    ,@('Tolak MBG, Ortu Siswa SD di Solo Pilih Bayar Rp 10 Ribu ke Dapur Sekolah', '2025-09-29T12:28:20+07:00', 'Tara Wahyu NV', 'https://news.detik.com/berita/d-8135453/tolak-mbg-ortu-siswa-sd-di-solo-pilih-bayar-rp-10-ribu-ke-dapur-sekolah')
    ,@('SPPG Sumberejo Semin Ditutup Sementara Buntut 19 Murid Diduga Keracunan MBG', '2025-09-29T12:22:19+07:00', 'Pradito Rida Pertana', 'https://www.detik.com/jogja/berita/d-8135426/sppg-sumberejo-semin-ditutup-sementara-buntut-19-murid-diduga-keracunan-mbg')
    ,@('Ortu Siswa SD Swasta di Solo Tolak MBG, Pilih Bayar Rp 10 Ribu ke Dapur Sekolah', '2025-09-29T12:00:38+07:00', 'Tara Wahyu NV', 'https://www.detik.com/jateng/berita/d-8135372/ortu-siswa-sd-swasta-di-solo-tolak-mbg-pilih-bayar-rp-10-ribu-ke-dapur-sekolah')
    ,@('Ternyata Ada Bakteri Salmonella di Sampel Makanan Program MBG di Bandung Barat', '2025-09-29T12:00:00+07:00', 'Antara, Fahri Zulfikar', 'https://www.detik.com/edu/sekolah/d-8135348/ternyata-ada-bakteri-salmonella-di-sampel-makanan-program-mbg-di-bandung-barat')
    ,@('Blak-blakan Menkop soal Tantangan Koperasi Merah Putih', '2025-09-26T17:03:11+07:00', 'Wisma Putra', 'https://www.detik.com/jabar/bisnis/d-8132102/blak-blakan-menkop-soal-tantangan-koperasi-merah-putih')
    ,@('Begini Strategi Khofifah Cegah Jeratan Pinjol Ilegal di Tengah Warga', '2025-09-26T10:45:16+07:00', 'Faiq Azmi', 'https://www.detik.com/jatim/berita/d-8131171/begini-strategi-khofifah-cegah-jeratan-pinjol-ilegal-di-tengah-warga')
    ,@('Kemenkop Puji Kopdes Aeng Batu-batu Takalar Raup Omzet Rp 2 Juta Per Hari', '2025-09-23T14:46:31+07:00', 'Adhe Junaedi Sholat', 'https://www.detik.com/sulsel/bisnis/d-8126210/kemenkop-puji-kopdes-aeng-batu-batu-takalar-raup-omzet-rp-2-juta-per-hari')
    ,@('Kemenkop Kawal Pemberian Pinjaman Modal 3.059 Koperasi Merah Putih di Sulsel', '2025-09-23T11:26:57+07:00', 'Adhe Junaedi Sholat', 'https://www.detik.com/sulsel/bisnis/d-8125752/kemenkop-kawal-pemberian-pinjaman-modal-3-059-koperasi-merah-putih-di-sulsel')
    ,@('Ujian Awal Koperasi Merah Putih di Usia Setipis Daun Kelor', '2025-09-22T10:30:00+07:00', 'Andry Haryanto', 'https://www.detik.com/jabar/berita/d-8123636/ujian-awal-koperasi-merah-putih-di-usia-setipis-daun-kelor')
    ,@('PPPK Bisa Diperbantukan di Kopdes, Pinjaman Rp 3 M Cair Pekan Depan', '2025-09-18T19:25:18+07:00', 'Agus Setyadi', 'https://www.detik.com/sumut/bisnis/d-8119042/pppk-bisa-diperbantukan-di-kopdes-pinjaman-rp-3-m-cair-pekan-depan')
    ,@('Melihat Koperasi Merah Putih di Melawai Jaksel', '2025-07-23T13:00:24+07:00', '-', 'https://finance.detik.com/foto-bisnis/d-8024832/melihat-koperasi-merah-putih-di-melawai-jaksel')
    ,@('Wamenkop Apresiasi 100% Pembentukan Badan Hukum Kopdes Merah Putih di Lahat', '2025-06-10T15:10:00+07:00', '-', 'https://www.detik.com/sumbagsel/foto/d-7957338/wamenkop-apresiasi-100-pembentukan-badan-hukum-kopdes-merah-putih-di-lahat')
    ,@('Momen Wamenkop Ferry Pantau Musdesus Kopdes Merah Putih di Padang', '2025-05-29T19:30:25+07:00', '-', 'https://www.detik.com/sumut/foto/d-7939188/momen-wamenkop-ferry-pantau-musdesus-kopdes-merah-putih-di-padang')
    ,@('Budi Arie Sambangi KPK, Minta Pengawalan Program Koperasi Desa', '2025-05-21T17:29:29+07:00', '-', 'https://news.detik.com/foto-news/d-7925409/budi-arie-sambangi-kpk-minta-pengawalan-program-koperasi-desa')
    ,@('20 Ribu Kopdes Ditargetkan Dapat Pinjaman dari Bank BUMN Tahun Ini', '2025-09-18T13:53:03+07:00', 'Retno Ayuningrum', 'https://finance.detik.com/moneter/d-8118291/20-ribu-kopdes-ditargetkan-dapat-pinjaman-dari-bank-bumn-tahun-ini')
    ,@('Link Pengumuman Hasil Seleksi Administrasi PMO Koperasi Merah Putih Kemenkop', '2025-09-15T14:30:00+07:00', 'Nikita Rosa', 'https://www.detik.com/edu/detikpedia/d-8112624/link-pengumuman-hasil-seleksi-administrasi-pmo-koperasi-merah-putih-kemenkop')
    ,@('Zulhas Minta Percepatan Pinjaman Kopdes Merah Putih: Rp 200 T Ada di Bank Himbara', '2025-09-15T12:18:25+07:00', 'Lisye Sri Rahayu', 'https://news.detik.com/berita/d-8112520/zulhas-minta-percepatan-pinjaman-kopdes-merah-putih-rp-200-t-ada-di-bank-himbara')
    ,@('Menteri Koperasi Tinjau Koperasi Merah Putih di Tuban, Ini Hasilnya', '2025-09-13T19:45:58+07:00', 'Ainur Rofiq', 'https://www.detik.com/jatim/bisnis/d-8110591/menteri-koperasi-tinjau-koperasi-merah-putih-di-tuban-ini-hasilnya')
    ,@('Kepala BGN Dipanggil Prabowo, Diminta Cegah Keracunan MBG dengan 4 Langkah Ini', '2025-09-29T09:30:24+07:00', 'Anggi Muliawati', 'https://www.detik.com/edu/sekolah/d-8135117/kepala-bgn-dipanggil-prabowo-diminta-cegah-keracunan-mbg-dengan-4-langkah-ini')
    ,@('Pidato Prabowo di PBB: Kemerdekaan Palestina dan Pesan Perdamaian', '2025-09-28T19:03:18+07:00', 'Sugiat Santoso', 'https://news.detik.com/kolom/d-8134771/pidato-prabowo-di-pbb-kemerdekaan-palestina-dan-pesan-perdamaian')
    ,@('Prabowo Tiba di Jakarta Usai Berpidato di Sidang Umum PBB', '2025-09-27T17:18:23+07:00', '-', 'https://news.detik.com/foto-news/d-8133403/prabowo-tiba-di-jakarta-usai-berpidato-di-sidang-umum-pbb')
    ,@('Potret Prabowo Didampingi Para Menteri Kabinet Merah Putih Hadiri Sidang Umum PBB Ke-80', '2025-09-24T08:25:14+07:00', '-', 'https://news.detik.com/foto-news/d-8127325/potret-prabowo-didampingi-para-menteri-kabinet-merah-putih-hadiri-sidang-umum-pbb-ke-80')
    ,@('Prabowo Serukan Keadilan Global dan Kesetaraan di Sidang Umum PBB ke-80', '2025-09-24T07:52:07+07:00', '-', 'https://news.detik.com/foto-news/d-8127294/prabowo-serukan-keadilan-global-dan-kesetaraan-di-sidang-umum-pbb-ke-80')
    ,@('Prabowo Bicara di PBB Soal Palestina, Usung Solusi Dua Negara', '2025-09-23T11:00:32+07:00', '-', 'https://news.detik.com/foto-news/d-8125619/prabowo-bicara-di-pbb-soal-palestina-usung-solusi-dua-negara')
    ,@('Potret Kunjungan Presiden Prabowo ke Area dan Korban Banjir Bali', '2025-09-14T06:06:36+07:00', '-', 'https://travel.detik.com/fototravel/d-8110851/potret-kunjungan-presiden-prabowo-ke-area-dan-korban-banjir-bali')
    ,@('Banyak Kasus Keracunan MBG di Sekolah, Begini Respons Presiden Prabowo', '2025-09-28T18:12:23+07:00', 'detikNews, Cicin Yulianti', 'https://www.detik.com/edu/sekolah/d-8134652/banyak-kasus-keracunan-mbg-di-sekolah-begini-respons-presiden-prabowo')
    ,@('Prabowo Bakal Panggil Kepala BGN Bahas Kasus Keracunan MBG', '2025-09-27T18:02:18+07:00', 'Eva Safitri', 'https://www.detik.com/sumut/berita/d-8133346/prabowo-bakal-panggil-kepala-bgn-bahas-kasus-keracunan-mbg')
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
}

Write-Host "Done writing rows"